$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted coin prices; some look numeric (e.g. "1.002"),
# so force text storage via NumberFormat, then restore the default style afterwards
# so no stray style index is left on the cell (matches source formatting).
$ws.Range('D2').Value = '30.224.43'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').Value = '1.896.24'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5181'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.58%  '
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08390'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.115'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.32'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +11.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.436'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.03%  '
$ws.Range('D14').Value = '1.893.43'
$ws.Range('E14').Value = '  -0.96%  '
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.28'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.13%  '
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06645'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.45%  '
$ws.Range('E20').Value = '  +1.32%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  -1.36%  '
$ws.Range('D23').Value = '30.217.04'
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.232'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('D26').Value = '2.110.16'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.351'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '129.76'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.093'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.23%  '
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.085'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.740'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02495'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06548'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.262'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('E38').Value = '  -0.24%  '
$ws.Range('E39').Value = '  -1.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.79'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.737'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6500'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.226'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6094'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.28'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.053'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.237'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.59'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.164'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.98'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.18%  '
